# feat(upload): subida a bbdd de los datos del excel
#
# The "id" column (column A) is no longer needed as a header in the
# template, since it will now be populated/managed by the database
# upload process. Remove the whole column so every following column
# shifts one place to the left (B->A, C->B, ... M->L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A ("id") entirely - this shifts nombre, descripcion,
# desc_sitio, horario, transporte, url, direccion, codpostal, latitud,
# longitud, fecha and autores one column to the left (A1:M1 -> A1:L1).
$ws.Range("A1").EntireColumn.Delete()
